$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the set_voltage (column G) values
$ws.Range("G3").Value = 54.4
$ws.Range("G4").Value = 54.4
$ws.Range("G19").Value = 51
$ws.Range("G24").Value = 54.4
$ws.Range("G25").Value = 54.4
$ws.Range("G26").Value = 54.4

# Update the active cell selection to I21
$ws.Range("I21").Select()
